# Se agregan los script 0125/0126/0127/0128/129/130 a la Clase Tests_AdmInstituciones
#
# The test data pool sheet keeps a "Verity1.0" password value in column C
# for the DEC_0078..DEC_0131 test rows. This edit bumps that shared value to
# "Verity1.1" for every one of those rows (rows 9-48), and it fills in the
# previously-empty D:J columns (BUSQUEDA..EMAIL_REPRESENTANTE) for the newly
# added script rows 41-48 with the default "SIN_DATO" placeholder, matching
# the pattern already used by rows 2-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password/version column for every existing Verity row
# (DEC_0078 .. DEC_0131) to the new "Verity1.1" value.
$ws.Range("C9:C48").Value = "Verity1.1"

# The newly added script rows (41-48, scripts DEC_0124..DEC_0131) only had
# A:C populated; fill the remaining D:J columns with the usual "SIN_DATO"
# placeholder used throughout the rest of the table.
$ws.Range("D41:J48").Value = "SIN_DATO"

# D41:J48 above already created these cells with the plain "SIN_DATO" text
# style used elsewhere in the column (style carried from the column
# default). The H41:H48 cells already existed (previously blank, styled
# like the other "quotePrefix" password cells), so copy the formatting
# used by the rest of the D:J columns onto them to keep styles consistent.
$ws.Range("D9").Copy()
$ws.Range("H41:H48").PasteSpecial(-4122)

# Reflect the editor's new scroll position/selection in the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("C50").Select()
